# Update the "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 2014
    4  = 124
    7  = 1675
    9  = 674
    12 = 25
    17 = 115
    18 = 135
    19 = 3884
    23 = 361
    24 = 722
    25 = 534
    27 = 33
    28 = 1678
    29 = 9
    30 = 26
    31 = 162
    32 = 12
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
